$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts B:M columns that
# currently hold the data one column to the right.
$ws.Columns("A:A").Insert()

# New header in A1
$ws.Range("A1").Value = "Date"

# New date values in A2:A4, formatted as dates (numFmtId 14 -> m/d/yy)
$ws.Range("A2").Value = (Get-Date -Year 2026 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A3").Value = (Get-Date -Year 2026 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A4").Value = (Get-Date -Year 2026 -Month 1 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("A2:A4").NumberFormat = "m/d/yy"

# Resize columns to fit content like the original bestFit columns
$ws.Columns("A:M").AutoFit()

# Update selection to match the target workbook view
$ws.Range("B11").Select()
